# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values for the second
# handoff/handback entry (row 3, the ef91d88f... file) on both the
# zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-18 03:41:30"
$wsZhCn.Range("G3").Value = "2016-01-18 03:42:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-18 03:41:47"
$wsDeDe.Range("G3").Value = "2016-01-18 03:42:59"
